# Insert a new price-record row before the existing row 149 (Poroto granado,
# Macroferia Regional de Talca). All rows from the old 149 through 211 shift
# down by one to 150..212, and the new row 149 carries the inserted record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 149:211 down to 150:212, leaving row 149 free for the new record.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the same fixed attributes used
# throughout this sheet (market, region, product, unit, origin, classification)
# plus the record-specific date / volume / price figures.
$ws.Range("A149").Value2 = 5
$ws.Range("B149").Value2 = "Macroferia Regional de Talca"
$ws.Range("C149").Value2 = "Maule"
$ws.Range("D149").Value2 = 45009
$ws.Range("E149").Value2 = 7
$ws.Range("F149").Value2 = 100112030
$ws.Range("G149").Value2 = "Poroto granado"
$ws.Range("H149").Value2 = "Sin especificar"
$ws.Range("I149").Value2 = "Primera"
$ws.Range("J149").Value2 = 300
$ws.Range("K149").Value2 = 28000
$ws.Range("L149").Value2 = 28000
$ws.Range("M149").Value2 = 28000
$ws.Range("N149").Value2 = "$/saco 25 kilos"
$ws.Range("O149").Value2 = "Región del Maule"
$ws.Range("P149").Value2 = 1120
$ws.Range("Q149").Value2 = 25
$ws.Range("R149").Value2 = "Hortaliza"
